$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.119053265332439
$ws.Range("C2").Value = 0.5699467346675675
$ws.Range("D2").Value = 0.920943356034627
$ws.Range("E2").Value = 1.160946734667561
$ws.Range("F2").Value = 0.6719467346675714
$ws.Range("G2").Value = 0.2209467346675638
$ws.Range("H2").Value = 0.3009467346675621
$ws.Range("B3").Value = 1.689000000000007
$ws.Range("C3").Value = 2.039996621367067
$ws.Range("D3").Value = 2.280000000000001
$ws.Range("E3").Value = 1.791000000000011
$ws.Range("F3").Value = 1.340000000000003
$ws.Range("G3").Value = 1.420000000000001
$ws.Range("B4").Value = 0.3509966213670595
$ws.Range("C4").Value = 0.590999999999994
$ws.Range("D4").Value = 0.1020000000000039
$ws.Range("E4").Value = -0.3490000000000037
$ws.Range("F4").Value = -0.2690000000000054
$ws.Range("G4").Value = -0.2489999999999952
$ws.Range("H4").Value = -0.1490000000000009
$ws.Range("I4").Value = -0.2489999999999952
$ws.Range("J4").Value = -0.4690015876295774
$ws.Range("B5").Value = 0.2400033786329345
$ws.Range("C5").Value = -0.2489966213670556
$ws.Range("D5").Value = -0.6999966213670632
$ws.Range("E5").Value = -0.6199966213670649
$ws.Range("F5").Value = -0.5999966213670547
$ws.Range("G5").Value = -0.4999966213670604
$ws.Range("H5").Value = -0.5999966213670547
$ws.Range("I5").Value = -0.8199982089966369
$ws.Range("B6").Value = -0.4889999999999901
$ws.Range("C6").Value = -0.9399999999999977
$ws.Range("D6").Value = -0.8599999999999994
$ws.Range("E6").Value = -0.8399999999999892
$ws.Range("F6").Value = -0.7399999999999949
$ws.Range("G6").Value = -0.8399999999999892
$ws.Range("H6").Value = -1.060001587629571
$ws.Range("B7").Value = -0.4510000000000076
$ws.Range("C7").Value = -0.3710000000000093
$ws.Range("D7").Value = -0.3509999999999991
$ws.Range("E7").Value = -0.2510000000000048
$ws.Range("F7").Value = -0.3509999999999991
$ws.Range("G7").Value = -0.5710015876295813
$ws.Range("B8").Value = 0.07999999999999828
$ws.Range("C8").Value = 0.1000000000000085
$ws.Range("D8").Value = 0.2000000000000028
$ws.Range("E8").Value = 0.1000000000000085
$ws.Range("F8").Value = -0.1200015876295737
$ws.Range("G8").Value = 0.100002285065301
$ws.Range("H8").Value = -0.000003441924292302279
$ws.Range("I8").Value = 0.1600000000000107
$ws.Range("B9").Value = 0.02000000000001019
$ws.Range("C9").Value = 0.1200000000000045
$ws.Range("D9").Value = 0.02000000000001019
$ws.Range("E9").Value = -0.200001587629572
$ws.Range("F9").Value = 0.02000228506530273
$ws.Range("G9").Value = -0.08000344192429058
$ws.Range("H9").Value = 0.08000000000001246
$ws.Range("B10").Value = 0.09999999999999432
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = -0.2200015876295822
$ws.Range("E10").Value = 0.000002285065292539912
$ws.Range("F10").Value = -0.1000034419243008
$ws.Range("G10").Value = 0.06000000000000227
$ws.Range("B11").Value = -0.09999999999999432
$ws.Range("C11").Value = -0.3200015876295765
$ws.Range("D11").Value = -0.09999771493470178
$ws.Range("E11").Value = -0.2000034419242951
$ws.Range("F11").Value = -0.03999999999999204
$ws.Range("B12").Value = -0.2200015876295822
$ws.Range("C12").Value = 0.000002285065292539912
$ws.Range("D12").Value = -0.1000034419243008
$ws.Range("E12").Value = 0.06000000000000227
$ws.Range("B13").Value = 0.2200038726948747
$ws.Range("C13").Value = 0.1199981457052814
$ws.Range("D13").Value = 0.2800015876295844
$ws.Range("B14").Value = -0.1000057269895933
$ws.Range("C14").Value = 0.05999771493470973
$ws.Range("B15").Value = 0.160003441924303
